# Actualizacion Datos Personales 4 nov
# Adds two new "rescatable" students to the "Rescatables" sheet, inserting
# them between the existing rows (new rows land at sheet rows 3 and 5),
# pushing the pre-existing students down to keep matricula order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert two blank rows to make room for the new students, keeping the
# existing three rows of data (currently rows 2-4) spread across rows 2-6.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(5).Insert()

# Row 2: existing student (unchanged matricula/grupo/materia)
$ws.Range("A2").Value = 21330051920003
$ws.Range("B2").Value = "BAEZ"
$ws.Range("C2").Value = "MARCELINO"
$ws.Range("D2").Value = "LUIS EDUARDO"
$ws.Range("E2").Value = "ÁLGEBRA"
$ws.Range("F2").Value = "1AV"
$ws.Range("G2").Value = 6

# Row 3: new student
$ws.Range("A3").Value = 21330051920007
$ws.Range("B3").Value = "COBOS"
$ws.Range("C3").Value = "NOLASCO"
$ws.Range("D3").Value = "YOLET"
$ws.Range("E3").Value = "ÁLGEBRA"
$ws.Range("F3").Value = "1AV"
$ws.Range("G3").Value = 6

# Row 4: existing student
$ws.Range("A4").Value = 21330051920017
$ws.Range("B4").Value = "MARTINEZ"
$ws.Range("C4").Value = "XOTLANIHUA"
$ws.Range("D4").Value = "YAIR"
$ws.Range("E4").Value = "ÁLGEBRA"
$ws.Range("F4").Value = "1AV"
$ws.Range("G4").Value = 6

# Row 5: new student
$ws.Range("A5").Value = 21330051920025
$ws.Range("B5").Value = "TAMAYO"
$ws.Range("C5").Value = "VARGAS"
$ws.Range("D5").Value = "JOSMAR JAHIR"
$ws.Range("E5").Value = "ÁLGEBRA"
$ws.Range("F5").Value = "1AV"
$ws.Range("G5").Value = 6

# Row 6: existing student
$ws.Range("A6").Value = 21330051920045
$ws.Range("B6").Value = "HERNANDEZ"
$ws.Range("C6").Value = "GALEOTE"
$ws.Range("D6").Value = "ZURY BETZABE"
$ws.Range("E6").Value = "ÁLGEBRA"
$ws.Range("F6").Value = "1BV"
$ws.Range("G6").Value = 6
